$wb = $excel.ActiveWorkbook

# --- Work on the "Transactions" sheet ---
$ws = $wb.Worksheets.Item("Transactions")
$ws.Activate() | Out-Null

# Update the transaction entries (Entry ID is now shown as a running
# sequence number instead of the raw loan-transaction id, and the
# Principal/Interest/Loan Balance figures were corrected to match the
# recalculated Summary sheet).
$ws.Range("A2").Value = 2
$ws.Range("F2").Value = 790.76
$ws.Range("G2").Value = 96.96
$ws.Range("J2").Value = 9209.24
$ws.Range("A3").Value = 1

# The stray formatted-but-empty cell in K2 is no longer needed.
$ws.Range("K2").Clear() | Out-Null

# Widen the columns so the longer transaction-date / loan-balance text fits.
$ws.Columns("C:C").AutoFit() | Out-Null
$ws.Columns("J:J").AutoFit() | Out-Null

# Leave the cursor on C2, as it ended up after the edits.
$ws.Range("C2").Select() | Out-Null
